# Delete column B ("Размер (точно как в ЛК, либо пусто)" / "37-38" sample)
# from the template sheet. Everything to the right (C:F) shifts one column
# to the left (becomes B:E), matching the committed template update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Delete()

# Leave the cursor where the author's session left it after the edit.
[void]$ws.Range("D11").Select()
